$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template gains three new leading columns (sample_id, project inserted
# after library_id, pushing source_well/conc_plate_* two columns to the
# right) plus three new trailing columns (bc_well, bc_seq, bc_read).

# Insert two blank columns at B:C - this shifts the existing source_well
# column (and its 96 well-id values) together with the two conc_plate_*
# columns two positions to the right, while preserving their original
# column widths/bestFit formatting automatically.
$ws.Range("B:C").Insert()

# Fill in the headers for the two newly inserted columns.
$ws.Range("B1").Value = "sample_id"
$ws.Range("C1").Value = "project"

# Add the three new trailing columns with their headers.
$ws.Range("G1").Value = "bc_well"
$ws.Range("H1").Value = "bc_seq"
$ws.Range("I1").Value = "bc_read"

# Widen the first three columns so the new headers are readable.
$ws.Range("A:A").ColumnWidth = 14.6
$ws.Range("B:B").ColumnWidth = 12.6
$ws.Range("C:C").ColumnWidth = 8.75

# Leave the selection where the author left it after adding the columns.
$ws.Range("G2").Select() | Out-Null
